# Auto-generated edit script: updates cryptos list per diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Plain text / non-numeric-looking values: direct assignment ---
$ws.Range('D2').Value = '66.905.37'
$ws.Range('E2').Value = '  -3.47%  '
$ws.Range('D3').Value = '3.527.78'
$ws.Range('E3').Value = '  -3.96%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -5.83%  '
$ws.Range('E6').Value = '  -3.18%  '
$ws.Range('D7').Value = '3.526.42'
$ws.Range('E7').Value = '  -3.95%  '
$ws.Range('E9').Value = '  -3.18%  '
$ws.Range('E10').Value = '  -3.23%  '
$ws.Range('E11').Value = '  -3.46%  '
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('E13').Value = '  -3.53%  '
$ws.Range('D14').Value = '4.133.60'
$ws.Range('E14').Value = '  -3.79%  '
$ws.Range('E15').Value = '  -2.24%  '
$ws.Range('D16').Value = '3.531.08'
$ws.Range('E16').Value = '  -4.14%  '
$ws.Range('D17').Value = '66.978.67'
$ws.Range('E17').Value = '  -3.42%  '
$ws.Range('E18').Value = '  +0.78%  '
$ws.Range('E19').Value = '  -2.66%  '
$ws.Range('E20').Value = '  -3.63%  '
$ws.Range('E21').Value = '  -4.61%  '
$ws.Range('E22').Value = '  -6.20%  '
$ws.Range('E23').Value = '  -2.80%  '
$ws.Range('E24').Value = '  -1.63%  '
$ws.Range('D25').Value = '3.676.59'
$ws.Range('E25').Value = '  -3.78%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('E27').Value = '  -0.92%  '
$ws.Range('E28').Value = '  -5.08%  '
$ws.Range('E29').Value = '  -8.80%  '
$ws.Range('E30').Value = '  -3.64%  '
$ws.Range('E31').Value = '  -2.20%  '
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('E33').Value = '  -4.47%  '
$ws.Range('B34').Value = 'Kaspa'
$ws.Range('C34').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('E34').Value = '  -3.62%  '
$ws.Range('B35').Value = 'ImmutableX'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('E35').Value = '  -5.89%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.529.25'
$ws.Range('E36').Value = '  -3.77%  '
$ws.Range('B37').Value = 'NEARProtocol'
$ws.Range('C37').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('E37').Value = '  -4.02%  '
$ws.Range('E38').Value = '  -4.76%  '
$ws.Range('E39').Value = '  +0.01%  '
$ws.Range('E40').Value = '  +0.06%  '
$ws.Range('E41').Value = '  -1.56%  '
$ws.Range('E42').Value = '  -3.06%  '
$ws.Range('E43').Value = '  -5.43%  '
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('E45').Value = '  -3.89%  '
$ws.Range('E46').Value = '  -2.09%  '
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('E48').Value = '  -3.34%  '
$ws.Range('E49').Value = '  -1.75%  '
$ws.Range('E50').Value = '  -3.30%  '
$ws.Range('E51').Value = '  -3.04%  '

# --- Numeric-looking text values: force text storage so Excel doesn't coerce them to numbers ---
$numericTextCells = @(
    'D4',
    'D5',
    'D6',
    'D9',
    'D10',
    'D11',
    'D12',
    'D13',
    'D15',
    'D19',
    'D20',
    'D21',
    'D22',
    'D23',
    'D24',
    'D27',
    'D28',
    'D29',
    'D30',
    'D31',
    'D33',
    'D34',
    'D35',
    'D37',
    'D38',
    'D40',
    'D41',
    'D42',
    'D43',
    'D44',
    'D45',
    'D46',
    'D47',
    'D48',
    'D49',
    'D50',
    'D51'
)
foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D4').Value = '1.00'
$ws.Range('D5').Value = '605.81'
$ws.Range('D6').Value = '153.15'
$ws.Range('D9').Value = '0.481'
$ws.Range('D10').Value = '0.140'
$ws.Range('D11').Value = '6.83'
$ws.Range('D12').Value = '0.426'
$ws.Range('D13').Value = '0.0000222'
$ws.Range('D15').Value = '31.64'
$ws.Range('D19').Value = '6.27'
$ws.Range('D20').Value = '15.30'
$ws.Range('D21').Value = '444.18'
$ws.Range('D22').Value = '9.25'
$ws.Range('D23').Value = '0.625'
$ws.Range('D24').Value = '78.01'
$ws.Range('D27').Value = '0.0000122'
$ws.Range('D28').Value = '10.16'
$ws.Range('D29').Value = '8.13'
$ws.Range('D30').Value = '2.53'
$ws.Range('D31').Value = '1.64'
$ws.Range('D33').Value = '25.58'
$ws.Range('D34').Value = '0.158'
$ws.Range('D35').Value = '1.87'
$ws.Range('D37').Value = '6.12'
$ws.Range('D38').Value = '7.96'
$ws.Range('D40').Value = '1.00'
$ws.Range('D41').Value = '175.01'
$ws.Range('D42').Value = '2.14'
$ws.Range('D43').Value = '5.52'
$ws.Range('D44').Value = '0.0859'
$ws.Range('D45').Value = '0.888'
$ws.Range('D46').Value = '45.67'
$ws.Range('D47').Value = '27.46'
$ws.Range('D48').Value = '2.59'
$ws.Range('D49').Value = '1.22'
$ws.Range('D50').Value = '7.53'
$ws.Range('D51').Value = '1.02'

foreach ($addr in $numericTextCells) {
    $ws.Range($addr).NumberFormat = "General"
    $ws.Range($addr).Style = "Normal"
}
